$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two input values that drive the whole calculator; all
# dependent formula cells (D6, D8, D9, E8, E9) recalc automatically.
$ws.Range("D3").Value = 666112.53
$ws.Range("D5").Value = 131681.51

# Move the live selection to match where the author left off editing.
$ws.Range("G7:I16").Select()
